# Add dev and prod profiles with H2 database configuration; enhance Swagger documentation
#
# Functionally (as far as this workbook is concerned) this updates the
# "Posts Report" sheet: two existing text cells get appended text, and two
# new data rows are appended to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: append text to the existing Title/Content values
$ws.Range("B2").Value = "фцвфцвawdawd"
$ws.Range("C2").Value = "фцвфвцww"

# New row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "awd"
$ws.Range("C5").Value = "awd"
$ws.Range("D5").Value = "2024-12-11T16:08:57.782325"

# New row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ыуаыуа"
$ws.Range("C6").Value = "ыуаыуа"
$ws.Range("D6").Value = "2024-12-11T16:21:41.180195"
